$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.212.35"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.283.72"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.59"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.99"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.81"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.23"
$ws.Range("E12").Value = "  +7.07%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.60"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "2.626.03"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.875"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "2.277.95"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "43.277.73"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  +4.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.58"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.52"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.76"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.70"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.90"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.95"
$ws.Range("E28").Value = "  -6.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.33"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.72"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.45"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0910"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  +5.76%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.69"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.94"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0359"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("E40").Value = "  +11.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.84"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.00"
$ws.Range("E42").Value = "  +5.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.241"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.51"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0998"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "1.430.30"
$ws.Range("E51").Value = "  +2.17%  "
